$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.045315207301049
$ws.Cells.Item(2, 4).Value = 1.051001048080566
$ws.Cells.Item(2, 5).Value = 1.048980915629172
$ws.Cells.Item(2, 6).Value = 1.05961281965858
$ws.Cells.Item(2, 9).Value = 1.039589985579217
$ws.Cells.Item(2, 10).Value = 1.050375598923911
$ws.Cells.Item(2, 11).Value = 1.053753762412953
$ws.Cells.Item(2, 12).Value = 1.051739247367915
$ws.Cells.Item(2, 13).Value = 1.062341847366663
$ws.Cells.Item(2, 14).Value = 1.020490181430688
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.046735611047428
$ws.Cells.Item(3, 4).Value = 1.052107221536946
$ws.Cells.Item(3, 5).Value = 1.050351359645765
$ws.Cells.Item(3, 6).Value = 1.060854367014548
$ws.Cells.Item(3, 9).Value = 1.039909136250726
$ws.Cells.Item(3, 10).Value = 1.051441289696072
$ws.Cells.Item(3, 11).Value = 1.054671739114369
$ws.Cells.Item(3, 12).Value = 1.05292040324517
$ws.Cells.Item(3, 13).Value = 1.063396578011904
$ws.Cells.Item(3, 14).Value = 1.020858155085652
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.047653769215566
$ws.Cells.Item(4, 4).Value = 1.052821971629605
$ws.Cells.Item(4, 5).Value = 1.051237556026769
$ws.Cells.Item(4, 6).Value = 1.06165697726501
$ws.Cells.Item(4, 9).Value = 1.04011384983022
$ws.Cells.Item(4, 10).Value = 1.052129501573666
$ws.Cells.Item(4, 11).Value = 1.055264126635324
$ws.Cells.Item(4, 12).Value = 1.05368358756658
$ws.Cells.Item(4, 13).Value = 1.06407774393469
$ws.Cells.Item(4, 14).Value = 1.021095452694806
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.048039543857213
$ws.Cells.Item(5, 4).Value = 1.053122212171867
$ws.Cells.Item(5, 5).Value = 1.051609980555573
$ws.Cells.Item(5, 6).Value = 1.061994217865025
$ws.Cells.Item(5, 9).Value = 1.040199482173262
$ws.Cells.Item(5, 10).Value = 1.052418503782546
$ws.Cells.Item(5, 11).Value = 1.055512785307245
$ws.Cells.Item(5, 12).Value = 1.05400417046951
$ws.Cells.Item(5, 13).Value = 1.064363794624899
$ws.Cells.Item(5, 14).Value = 1.021195020839869
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.048104304431597
$ws.Cells.Item(6, 4).Value = 1.053172609879753
$ws.Cells.Item(6, 5).Value = 1.051672504641348
$ws.Cells.Item(6, 6).Value = 1.062050831789324
$ws.Cells.Item(6, 9).Value = 1.040213835085694
$ws.Cells.Item(6, 10).Value = 1.052467009750235
$ws.Cells.Item(6, 11).Value = 1.055554513939095
$ws.Cells.Item(6, 12).Value = 1.054057982636073
$ws.Cells.Item(6, 13).Value = 1.064411805629424
$ws.Cells.Item(6, 14).Value = 1.021211727558128
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.047658924809086
$ws.Cells.Item(7, 4).Value = 1.052825984398887
$ws.Cells.Item(7, 5).Value = 1.051242532896154
$ws.Cells.Item(7, 6).Value = 1.061661484180284
$ws.Cells.Item(7, 9).Value = 1.040114995738869
$ws.Cells.Item(7, 10).Value = 1.052133364496066
$ws.Cells.Item(7, 11).Value = 1.05526745071801
$ws.Cells.Item(7, 12).Value = 1.053687872226645
$ws.Cells.Item(7, 13).Value = 1.064081567377609
$ws.Cells.Item(7, 14).Value = 1.021096783881301
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.045795436949809
$ws.Cells.Item(8, 4).Value = 1.051375097398992
$ws.Cells.Item(8, 5).Value = 1.049444184945166
$ws.Cells.Item(8, 6).Value = 1.060032564757799
$ws.Cells.Item(8, 9).Value = 1.039698217141147
$ws.Cells.Item(8, 10).Value = 1.050736038293257
$ws.Cells.Item(8, 11).Value = 1.054064331118757
$ws.Cells.Item(8, 12).Value = 1.052138655317613
$ws.Cells.Item(8, 13).Value = 1.062698572944323
$ws.Cells.Item(8, 14).Value = 1.020614707438209
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.042504301441656
$ws.Cells.Item(9, 4).Value = 1.048810501344853
$ws.Cells.Item(9, 5).Value = 1.046270672817325
$ws.Cells.Item(9, 6).Value = 1.057156258367839
$ws.Cells.Item(9, 9).Value = 1.038949970834522
$ws.Cells.Item(9, 10).Value = 1.048263172475871
$ws.Cells.Item(9, 11).Value = 1.051931848215661
$ws.Cells.Item(9, 12).Value = 1.049400106367906
$ws.Cells.Item(9, 13).Value = 1.060251331825359
$ws.Cells.Item(9, 14).Value = 1.019759001468398
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.040304866432156
$ws.Cells.Item(10, 4).Value = 1.047095214316144
$ws.Cells.Item(10, 5).Value = 1.04415161641191
$ws.Cells.Item(10, 6).Value = 1.055234487863284
$ws.Cells.Item(10, 9).Value = 1.038441761398143
$ws.Cells.Item(10, 10).Value = 1.046607226235192
$ws.Cells.Item(10, 11).Value = 1.05050163380803
$ws.Cells.Item(10, 12).Value = 1.047568340827816
$ws.Cells.Item(10, 13).Value = 1.058612742742771
$ws.Cells.Item(10, 14).Value = 1.019184269769427
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.039351141146194
$ws.Cells.Item(11, 4).Value = 1.046351110127944
$ws.Cells.Item(11, 5).Value = 1.043233173582061
$ws.Cells.Item(11, 6).Value = 1.054401282379497
$ws.Cells.Item(11, 9).Value = 1.038219458289838
$ws.Cells.Item(11, 10).Value = 1.045888380665831
$ws.Cells.Item(11, 11).Value = 1.049880259841437
$ws.Cells.Item(11, 12).Value = 1.046773668062659
$ws.Cells.Item(11, 13).Value = 1.057901483495754
$ws.Cells.Item(11, 14).Value = 1.018934376858318
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.038996674679487
$ws.Cells.Item(12, 4).Value = 1.046074506097192
$ws.Cells.Item(12, 5).Value = 1.042891885434782
$ws.Cells.Item(12, 6).Value = 1.054091627866523
$ws.Cells.Item(12, 9).Value = 1.038136546077181
$ws.Cells.Item(12, 10).Value = 1.045621092692705
$ws.Cells.Item(12, 11).Value = 1.049649137144625
$ws.Cells.Item(12, 12).Value = 1.046478259736094
$ws.Cells.Item(12, 13).Value = 1.057637024485949
$ws.Cells.Item(12, 14).Value = 1.018841399198962
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.039072718571844
$ws.Cells.Item(13, 4).Value = 1.046133848199396
$ws.Cells.Item(13, 5).Value = 1.042965099219924
$ws.Cells.Item(13, 6).Value = 1.05415805734877
$ws.Cells.Item(13, 9).Value = 1.038154346392529
$ws.Cells.Item(13, 10).Value = 1.045678439469965
$ws.Cells.Item(13, 11).Value = 1.049698728145097
$ws.Cells.Item(13, 12).Value = 1.046541636383512
$ws.Cells.Item(13, 13).Value = 1.057693763932763
$ws.Cells.Item(13, 14).Value = 1.018861350320142
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.039321845169657
$ws.Cells.Item(14, 4).Value = 1.046328250273357
$ws.Cells.Item(14, 5).Value = 1.043204965437535
$ws.Cells.Item(14, 6).Value = 1.054375689632125
$ws.Cells.Item(14, 9).Value = 1.038212611664304
$ws.Cells.Item(14, 10).Value = 1.045866292231526
$ws.Cells.Item(14, 11).Value = 1.049861161651768
$ws.Cells.Item(14, 12).Value = 1.046749254270465
$ws.Cells.Item(14, 13).Value = 1.057879628676015
$ws.Cells.Item(14, 14).Value = 1.018926694493389
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.03947531220068
$ws.Cells.Item(15, 4).Value = 1.046447999772819
$ws.Cells.Item(15, 5).Value = 1.043352736526189
$ws.Cells.Item(15, 6).Value = 1.05450975807244
$ws.Cells.Item(15, 9).Value = 1.038248465858995
$ws.Cells.Item(15, 10).Value = 1.045981997705303
$ws.Cells.Item(15, 11).Value = 1.049961200205438
$ws.Cells.Item(15, 12).Value = 1.046877143673879
$ws.Cells.Item(15, 13).Value = 1.057994110746299
$ws.Cells.Item(15, 14).Value = 1.01896693443759
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.040368132902271
$ws.Cells.Item(16, 4).Value = 1.047144568759532
$ws.Cells.Item(16, 5).Value = 1.044212551460283
$ws.Cells.Item(16, 6).Value = 1.055289762102191
$ws.Cells.Item(16, 9).Value = 1.038456467482235
$ws.Cells.Item(16, 10).Value = 1.04665489510607
$ws.Cells.Item(16, 11).Value = 1.050542828153117
$ws.Cells.Item(16, 12).Value = 1.047621048465985
$ws.Cells.Item(16, 13).Value = 1.058659909660651
$ws.Cells.Item(16, 14).Value = 1.019200832502733
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.040927808072583
$ws.Cells.Item(17, 4).Value = 1.047581137350863
$ws.Cells.Item(17, 5).Value = 1.044751651683255
$ws.Cells.Item(17, 6).Value = 1.055778749257074
$ws.Cells.Item(17, 9).Value = 1.038586339179443
$ws.Cells.Item(17, 10).Value = 1.047076497920138
$ws.Cells.Item(17, 11).Value = 1.05090710786319
$ws.Cells.Item(17, 12).Value = 1.048087273382549
$ws.Cells.Item(17, 13).Value = 1.059077078999969
$ws.Cells.Item(17, 14).Value = 1.019347273674225
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.041254126814615
$ws.Cells.Item(18, 4).Value = 1.04783564818169
$ws.Cells.Item(18, 5).Value = 1.045066015455941
$ws.Cells.Item(18, 6).Value = 1.056063864642474
$ws.Cells.Item(18, 9).Value = 1.038661874638668
$ws.Cells.Item(18, 10).Value = 1.047322237404355
$ws.Cells.Item(18, 11).Value = 1.05111938547364
$ws.Cells.Item(18, 12).Value = 1.048359069515197
$ws.Cells.Item(18, 13).Value = 1.059320239289861
$ws.Cells.Item(18, 14).Value = 1.019432591027006
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.041365371173623
$ws.Cells.Item(19, 4).Value = 1.047922407494649
$ws.Cells.Item(19, 5).Value = 1.045173191258932
$ws.Cells.Item(19, 6).Value = 1.05616106436859
$ws.Cells.Item(19, 9).Value = 1.038687593604722
$ws.Cells.Item(19, 10).Value = 1.047405998797454
$ws.Cells.Item(19, 11).Value = 1.051191732750919
$ws.Cells.Item(19, 12).Value = 1.048451720538761
$ws.Cells.Item(19, 13).Value = 1.059403122379684
$ws.Cells.Item(19, 14).Value = 1.0194616652412
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.040867773752061
$ws.Cells.Item(20, 4).Value = 1.04753431138194
$ws.Cells.Item(20, 5).Value = 1.044693820089301
$ws.Cells.Item(20, 6).Value = 1.05572629620453
$ws.Cells.Item(20, 9).Value = 1.03857242758147
$ws.Cells.Item(20, 10).Value = 1.047031281992213
$ws.Cells.Item(20, 11).Value = 1.050868044873881
$ws.Cells.Item(20, 12).Value = 1.048037266894298
$ws.Cells.Item(20, 13).Value = 1.059032338037644
$ws.Cells.Item(20, 14).Value = 1.019331572205037
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.039248489459616
$ws.Cells.Item(21, 4).Value = 1.046271009555812
$ws.Cells.Item(21, 5).Value = 1.043134334659491
$ws.Cells.Item(21, 6).Value = 1.054311606939054
$ws.Cells.Item(21, 9).Value = 1.038195463363661
$ws.Cells.Item(21, 10).Value = 1.045810981952957
$ws.Cells.Item(21, 11).Value = 1.049813337808317
$ws.Cells.Item(21, 12).Value = 1.046688122409287
$ws.Cells.Item(21, 13).Value = 1.05782490351009
$ws.Cells.Item(21, 14).Value = 1.018907456590063
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.038229160272466
$ws.Cells.Item(22, 4).Value = 1.045475500946044
$ws.Cells.Item(22, 5).Value = 1.042153025359034
$ws.Cells.Item(22, 6).Value = 1.053421179240823
$ws.Cells.Item(22, 9).Value = 1.037956489171895
$ws.Cells.Item(22, 10).Value = 1.045042127455022
$ws.Cells.Item(22, 11).Value = 1.049148367065447
$ws.Cells.Item(22, 12).Value = 1.045838519864797
$ws.Cells.Item(22, 13).Value = 1.057064202322182
$ws.Cells.Item(22, 14).Value = 1.018639893173933
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.038769643964927
$ws.Cells.Item(23, 4).Value = 1.045897332233693
$ws.Cells.Item(23, 5).Value = 1.042673313647043
$ws.Cells.Item(23, 6).Value = 1.053893304050173
$ws.Cells.Item(23, 9).Value = 1.038083360422627
$ws.Cells.Item(23, 10).Value = 1.04544986536884
$ws.Cells.Item(23, 11).Value = 1.049501055870193
$ws.Cells.Item(23, 12).Value = 1.046289039181876
$ws.Cells.Item(23, 13).Value = 1.057467611766221
$ws.Cells.Item(23, 14).Value = 1.018781819904806
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.040894901080364
$ws.Cells.Item(24, 4).Value = 1.047555470430673
$ws.Cells.Item(24, 5).Value = 1.044719951958028
$ws.Cells.Item(24, 6).Value = 1.055749997798416
$ws.Cells.Item(24, 9).Value = 1.038578714302898
$ws.Cells.Item(24, 10).Value = 1.04705171366264
$ws.Cells.Item(24, 11).Value = 1.05088569637738
$ws.Cells.Item(24, 12).Value = 1.048059863122124
$ws.Cells.Item(24, 13).Value = 1.059052555070649
$ws.Cells.Item(24, 14).Value = 1.019338667330078
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.043356055272752
$ws.Cells.Item(25, 4).Value = 1.049474474555147
$ws.Cells.Item(25, 5).Value = 1.047091676197783
$ws.Cells.Item(25, 6).Value = 1.057900581073501
$ws.Cells.Item(25, 9).Value = 1.039145057529566
$ws.Cells.Item(25, 10).Value = 1.048903748486816
$ws.Cells.Item(25, 11).Value = 1.052484640004114
$ws.Cells.Item(25, 12).Value = 1.05010913672689
$ws.Cells.Item(25, 13).Value = 1.060885236777947
$ws.Cells.Item(25, 14).Value = 1.019980967459269
